$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 holds a sequence of "blog" entries referencing serial numbers.
# A new blog post (ser: 155) is inserted at the head (B11), shifting the
# existing entries back one slot:
#   B11 (was "ser: 153") -> becomes the new head entry "ser: 155"
#   D11 (was "ser: 152") -> takes over B11's old value "ser: 153"
#   I11 (was "ser: 151") -> takes over D11's old value "ser: 152"
$ws.Range("B11").Value = "type: blog`nwidth: 2`nheight: 1`nser: 155"
$ws.Range("D11").Value = "type: blog`nwidth: 2`nheight: 1`nser: 153"
$ws.Range("I11").Value = "type: blog`nwidth: 2`nheight: 1`nser: 152"
